$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Cd2"
$ws.Range("C2").Value = "Cd48"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7139163333333333
$ws.Range("H2").Value = 2.141749
$ws.Range("I2").Value = 0.5325697612629201
$ws.Range("J2").Value = 0.5325697612629202
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06575866666666667
$ws.Range("N2").Value = 0.197276
$ws.Range("O2").Value = 0.0003210775944521127
$ws.Range("P2").Value = 0.000321193332710493
$ws.Range("Q2").Value = 0.04694618619155556
$ws.Range("R2").Value = 0.422515675724
$ws.Range("S2").Value = 0.0001709962178242343
$ws.Range("T2").Value = 0.0001710578565208689

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Cd2"
$ws.Range("C3").Value = "Cd48"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7139163333333333
$ws.Range("H3").Value = 2.141749
$ws.Range("I3").Value = 0.5325697612629201
$ws.Range("J3").Value = 0.5325697612629202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 93.51811733333334
$ws.Range("N3").Value = 280.554352
$ws.Range("O3").Value = 0.4566177155519742
$ws.Range("P3").Value = 0.4567823117120824
$ws.Range("Q3").Value = 66.76411142684978
$ws.Range("R3").Value = 600.8770028416479
$ws.Range("S3").Value = 0.2431807877599348
$ws.Range("T3").Value = 0.2432684466976285

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Cd2"
$ws.Range("C4").Value = "Cd48"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7139163333333333
$ws.Range("H4").Value = 2.141749
$ws.Range("I4").Value = 0.5325697612629201
$ws.Range("J4").Value = 0.5325697612629202
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 111.0008776666667
$ws.Range("N4").Value = 333.002633
$ws.Range("O4").Value = 0.5419801919638461
$ws.Range("P4").Value = 0.5421755585810701
$ws.Range("Q4").Value = 79.24533958056855
$ws.Range("R4").Value = 713.2080562251169
$ws.Range("S4").Value = 0.2886422614434171
$ws.Range("T4").Value = 0.2887463077961109

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Cd2"
$ws.Range("C5").Value = "Cd48"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7139163333333333
$ws.Range("H5").Value = 2.141749
$ws.Range("I5").Value = 0.5325697612629201
$ws.Range("J5").Value = 0.5325697612629202
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.2213985
$ws.Range("N5").Value = 0.442797
$ws.Range("O5").Value = 0.001081014889727682
$ws.Range("P5").Value = 0.0007209363741367839
$ws.Range("Q5").Value = 0.1580600053255
$ws.Range("R5").Value = 0.9483600319529999
$ws.Range("S5").Value = 0.0005757158417439334
$ws.Range("T5").Value = 0.0003839489126597823

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cd2"
$ws.Range("C6").Value = "Cd48"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6265959999999999
$ws.Range("H6").Value = 1.879788
$ws.Range("I6").Value = 0.4674302387370798
$ws.Range("J6").Value = 0.4674302387370799
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06575866666666667
$ws.Range("N6").Value = 0.197276
$ws.Range("O6").Value = 0.0003210775944521127
$ws.Range("P6").Value = 0.000321193332710493
$ws.Range("Q6").Value = 0.04120411749866667
$ws.Range("R6").Value = 0.370837057488
$ws.Range("S6").Value = 0.0001500813766278783
$ws.Range("T6").Value = 0.0001501354761896241

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd2"
$ws.Range("C7").Value = "Cd48"
$ws.Range("D7").Value = "M1"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6265959999999999
$ws.Range("H7").Value = 1.879788
$ws.Range("I7").Value = 0.4674302387370798
$ws.Range("J7").Value = 0.4674302387370799
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 93.51811733333334
$ws.Range("N7").Value = 280.554352
$ws.Range("O7").Value = 0.4566177155519742
$ws.Range("P7").Value = 0.4567823117120824
$ws.Range("Q7").Value = 58.59807824859733
$ws.Range("R7").Value = 527.382704237376
$ws.Range("S7").Value = 0.2134369277920393
$ws.Range("T7").Value = 0.2135138650144539

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd2"
$ws.Range("C8").Value = "Cd48"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6265959999999999
$ws.Range("H8").Value = 1.879788
$ws.Range("I8").Value = 0.4674302387370798
$ws.Range("J8").Value = 0.4674302387370799
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 111.0008776666667
$ws.Range("N8").Value = 333.002633
$ws.Range("O8").Value = 0.5419801919638461
$ws.Range("P8").Value = 0.5421755585810701
$ws.Range("Q8").Value = 69.55270594242266
$ws.Range("R8").Value = 625.974353481804
$ws.Range("S8").Value = 0.2533379305204289
$ws.Range("T8").Value = 0.2534292507849593

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd2"
$ws.Range("C9").Value = "Cd48"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6265959999999999
$ws.Range("H9").Value = 1.879788
$ws.Range("I9").Value = 0.4674302387370798
$ws.Range("J9").Value = 0.4674302387370799
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.2213985
$ws.Range("N9").Value = 0.442797
$ws.Range("O9").Value = 0.001081014889727682
$ws.Range("P9").Value = 0.0007209363741367839
$ws.Range("Q9").Value = 0.138727414506
$ws.Range("R9").Value = 0.8323644870359999
$ws.Range("S9").Value = 0.0005052990479837483
$ws.Range("T9").Value = 0.0003369874614770016
